$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated DM_Stat (C) and P_Value (D) values per row
$ws.Range("C2").Value = 0.4670143626189159
$ws.Range("D2").Value = 0.6450817842442591

$ws.Range("C3").Value = 0.6273349441042665
$ws.Range("D3").Value = 0.5368969243698369

$ws.Range("C4").Value = 1.125520101214776
$ws.Range("D4").Value = 0.2724980123814509

$ws.Range("C5").Value = 0.2985078904257442
$ws.Range("D5").Value = 0.7681168213801954

$ws.Range("C6").Value = 0.1996948759527565
$ws.Range("D6").Value = 0.8435535027428456

$ws.Range("C7").Value = 0.6670507879741072
$ws.Range("D7").Value = 0.5116760357180099

$ws.Range("C8").Value = -0.1473036109708511
$ws.Range("D8").Value = 0.8842340124100734

$ws.Range("C9").Value = 0.3935165975158829
$ws.Range("D9").Value = 0.6977257961456327

$ws.Range("C10").Value = -0.2884709506476925
$ws.Range("D10").Value = 0.7756867186111696

$ws.Range("C11").Value = -0.7154783907819404
$ws.Range("D11").Value = 0.4818386387408213

$wb.Save()
